$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value is a plain number would be auto-converted
# from Text to Number by Excels normal type inference. The source workbook
# stores every Price/Volume cell as text, so for those specific cells we briefly
# force a Text number format, assign the literal string, then restore General -
# mirroring exactly what typing the same text in real Excel (with a leading
# apostrophe) would do.

$ws.Range("D2").Value = "29.350.20"
$ws.Range("D3").Value = "1.843.20"
$ws.Range("E3").Value = "  -0.73%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9987"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  -0.51%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.05"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -0.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6297"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -0.62%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9998"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07442"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -1.41%  "
$ws.Range("E9").Value = "  -0.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.84"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +1.07%  "
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").Value = "1.847.60"
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.983"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6787"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -0.95%  "
$ws.Range("E15").Value = "  -2.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.05"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -1.69%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.259"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  +1.76%  "
$ws.Range("D18").Value = "29.325.73"
$ws.Range("E18").Value = "  -0.74%  "
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.31"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -0.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9998"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -0.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.401"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -1.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -0.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "158.08"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -0.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.482"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1352"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -2.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.41"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -1.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06542"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +14.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.445"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +1.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.487"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +0.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.065"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.056"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -0.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.838"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +0.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.139"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -1.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6966"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("E36").Value = "  -0.97%  "
$ws.Range("E37").Value = "  +1.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.815"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -0.31%  "
$ws.Range("D39").Value = "1.243.77"
$ws.Range("E39").Value = "  -0.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.798"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +4.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9305"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +2.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9994"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -0.39%  "
$ws.Range("D43").Value = "1.990.00"
$ws.Range("E43").Value = "  -1.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.73"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -0.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.58"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -0.97%  "
$ws.Range("E46").Value = "  +4.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.047"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -1.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.711"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +1.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.021"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -0.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1145"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -2.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3891"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -2.08%  "
